$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ D = "..."; E = "..." }  (D is optional per row)
$updates = @{
    2  = @{ D = "42.547.21";  E = "  -0.39%  " }
    3  = @{ D = "2.295.58";   E = "  -0.49%  " }
    4  = @{                   E = "  +0.04%  " }
    5  = @{ D = "316.71";     E = "  -0.62%  " }
    6  = @{ D = "104.35";     E = "  +0.28%  " }
    7  = @{ D = "0.628";      E = "  -0.70%  " }
    8  = @{                   E = "  -0.02%  " }
    9  = @{                   E = "  -1.48%  " }
    10 = @{ D = "39.68";      E = "  -0.37%  " }
    11 = @{                   E = "  -0.48%  " }
    12 = @{                   E = "  +1.69%  " }
    13 = @{                   E = "  +2.55%  " }
    14 = @{                   E = "  +3.49%  " }
    15 = @{ D = "15.39";      E = "  +0.06%  " }
    16 = @{ D = "2.643.85";   E = "  -0.55%  " }
    17 = @{ D = "2.292.21";   E = "  -0.64%  " }
    18 = @{ D = "42.645.92";  E = "  -0.17%  " }
    19 = @{ D = "14.89";      E = "  +34.78%  " }
    20 = @{                   E = "  -0.15%  " }
    21 = @{                   E = "  +0.01%  " }
    22 = @{ D = "73.94";      E = "  +0.93%  " }
    23 = @{ D = "3.54";       E = "  -1.58%  " }
    24 = @{ D = "266.01";     E = "  -5.72%  " }
    25 = @{                   E = "  -2.82%  " }
    26 = @{ D = "1.01";       E = "  +0.30%  " }
    27 = @{ D = "10.94";      E = "  +0.32%  " }
    28 = @{                   E = "  -0.19%  " }
    29 = @{ D = "6.81";       E = "  +14.99%  " }
    30 = @{ D = "22.59";      E = "  -1.58%  " }
    31 = @{                   E = "  +3.04%  " }
    32 = @{ D = "166.25";     E = "  +0.68%  " }
    33 = @{                   E = "  +0.20%  " }
    34 = @{                   E = "  -4.79%  " }
    35 = @{                   E = "  -0.23%  " }
    36 = @{ D = "0.114";      E = "  -3.01%  " }
    37 = @{                   E = "  -2.19%  " }
    38 = @{                   E = "  -5.31%  " }
    39 = @{ D = "3.72";       E = "  -0.41%  " }
    40 = @{ D = "2.69";       E = "  -3.10%  " }
    41 = @{                   E = "  +3.21%  " }
    42 = @{ D = "70.42";      E = "  +0.58%  " }
    43 = @{                   E = "  +0.52%  " }
    44 = @{ D = "95.24";      E = "  -3.24%  " }
    45 = @{                   E = "  -0.05%  " }
    46 = @{ D = "12.31";      E = "  +1.07%  " }
    47 = @{ D = "114.90";     E = "  +2.11%  " }
    48 = @{ D = "79.78";      E = "  -0.14%  " }
    49 = @{ D = "1.702.55";   E = "  +5.82%  " }
    50 = @{                   E = "  -1.87%  " }
    51 = @{ D = "5.10";       E = "  -4.08%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]

    if ($vals.ContainsKey("D")) {
        $dCell = $ws.Range("D$row")
        # Force text storage so numeric-looking strings (e.g. "5.10", "316.71")
        # keep their exact original formatting instead of being parsed as numbers.
        $dCell.NumberFormat = "@"
        $dCell.Value = $vals["D"]
        # Restore the cell's original (default) style now that the text value
        # is locked in, so no stray style/number-format survives on the cell.
        $dCell.Style = "Normal"
    }

    if ($vals.ContainsKey("E")) {
        $ws.Range("E$row").Value = $vals["E"]
    }
}
